$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1688.8667
$ws.Range("I19").Value = 1126.3334
$ws.Range("J19").Value = 2063.889
$ws.Range("K19").Value = 1126.3334
$ws.Range("L19").Value = 2063.889
$ws.Range("M19").Value = -951.3334
$ws.Range("N19").Value = -2413.889

$ws.Range("H29").Value = 17987.5
$ws.Range("I29").Value = 4444
$ws.Range("J29").Value = 21373.375
$ws.Range("K29").Value = 13332
$ws.Range("L29").Value = 64120.125
$ws.Range("M29").Value = -13051
$ws.Range("N29").Value = -64682.125

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H80").Value = 2241872.8
$ws.Range("I80").Value = 1634517.6
$ws.Range("J80").Value = 2646776.2
$ws.Range("K80").Value = 4903552.800000001
$ws.Range("L80").Value = 7940328.600000001
$ws.Range("M80").Value = -4902554.800000001
$ws.Range("N80").Value = -7942324.600000001

$ws.Range("H83").Value = 2241872.8
$ws.Range("I83").Value = 1634517.6
$ws.Range("J83").Value = 2646776.2
$ws.Range("K83").Value = 14710658.4
$ws.Range("L83").Value = 23820985.8
$ws.Range("M83").Value = -14705666.4
$ws.Range("N83").Value = -23830969.8

$ws.Range("H86").Value = 8382.666999999999
$ws.Range("I86").Value = 4139.8
$ws.Range("J86").Value = 11413.286
$ws.Range("K86").Value = 4139.8
$ws.Range("L86").Value = 11413.286
$ws.Range("M86").Value = -3016.8
$ws.Range("N86").Value = -13659.286

$ws.Range("H88").Value = 2589.9375
$ws.Range("I88").Value = 2257.375
$ws.Range("J88").Value = 2922.5
$ws.Range("K88").Value = 2257.375
$ws.Range("L88").Value = 2922.5
$ws.Range("M88").Value = -1851.375
$ws.Range("N88").Value = -3734.5

$ws.Range("H89").Value = 8382.666999999999
$ws.Range("I89").Value = 4139.8
$ws.Range("J89").Value = 11413.286
$ws.Range("K89").Value = 20699
$ws.Range("L89").Value = 57066.43
$ws.Range("M89").Value = -15083
$ws.Range("N89").Value = -68298.42999999999

$ws.Range("H91").Value = 2589.9375
$ws.Range("I91").Value = 2257.375
$ws.Range("J91").Value = 2922.5
$ws.Range("K91").Value = 2257.375
$ws.Range("L91").Value = 2922.5
$ws.Range("M91").Value = -853.375
$ws.Range("N91").Value = -5730.5

$ws.Range("H132").Value = 7179.3687
$ws.Range("I132").Value = 3959.4167
$ws.Range("K132").Value = 11878.2501
$ws.Range("M132").Value = -9348.250100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3587.145
$ws.Range("I32").Value = 3151.5078
$ws.Range("K32").Value = 3151.5078
$ws.Range("M32").Value = -2864.5078

$ws.Range("H43").Value = 37504.332
$ws.Range("J43").Value = 37021.285
$ws.Range("L43").Value = 37021.285
$ws.Range("N43").Value = -37647.285

$ws.Range("H61").Value = 10478488
$ws.Range("I61").Value = 13127529
$ws.Range("J61").Value = 2001558.6
$ws.Range("K61").Value = 13127529
$ws.Range("L61").Value = 2001558.6
$ws.Range("M61").Value = -13127317
$ws.Range("N61").Value = -2001982.6

$ws.Range("H63").Value = 4161.625
$ws.Range("I63").Value = 4085.1428
$ws.Range("J63").Value = 4697
$ws.Range("K63").Value = 4085.1428
$ws.Range("L63").Value = 4697
$ws.Range("M63").Value = -3399.1428
$ws.Range("N63").Value = -6069

$ws.Range("H66").Value = 4161.625
$ws.Range("I66").Value = 4085.1428
$ws.Range("J66").Value = 4697
$ws.Range("K66").Value = 20425.714
$ws.Range("L66").Value = 23485
$ws.Range("M66").Value = -16993.714
$ws.Range("N66").Value = -30349

$ws.Range("H88").Value = 3132.2856
$ws.Range("I88").Value = 1799
$ws.Range("J88").Value = 3665.6
$ws.Range("K88").Value = 1799
$ws.Range("L88").Value = 3665.6
$ws.Range("M88").Value = -1393
$ws.Range("N88").Value = -4477.6

$ws.Range("H91").Value = 3132.2856
$ws.Range("I91").Value = 1799
$ws.Range("J91").Value = 3665.6
$ws.Range("K91").Value = 1799
$ws.Range("L91").Value = 3665.6
$ws.Range("M91").Value = -395
$ws.Range("N91").Value = -6473.6

$ws.Range("H97").Value = 1856.5217
$ws.Range("I97").Value = 1359.4445
$ws.Range("K97").Value = 1359.4445
$ws.Range("M97").Value = -863.4445000000001

$ws.Range("H110").Value = 1571.8125
$ws.Range("I110").Value = 512.6667
$ws.Range("K110").Value = 512.6667
$ws.Range("M110").Value = 1532.3333

$ws.Range("H123").Value = 325000
$ws.Range("J123").Value = 325000
$ws.Range("L123").Value = 325000
$ws.Range("N123").Value = -334800

$ws.Range("H136").Value = 10478488
$ws.Range("I136").Value = 13127529
$ws.Range("J136").Value = 2001558.6
$ws.Range("K136").Value = 39382587
$ws.Range("L136").Value = 6004675.800000001
$ws.Range("M136").Value = -39380037
$ws.Range("N136").Value = -6009775.800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5324.1904
$ws.Range("I20").Value = 6860.7334
$ws.Range("K20").Value = 6860.7334
$ws.Range("M20").Value = -6613.7334

$ws.Range("H86").Value = 43427.47
$ws.Range("I86").Value = 62608.727
$ws.Range("K86").Value = 62608.727
$ws.Range("M86").Value = -61485.727

$ws.Range("H89").Value = 43427.47
$ws.Range("I89").Value = 62608.727
$ws.Range("K89").Value = 313043.635
$ws.Range("M89").Value = -307427.635

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27030068
$ws.Range("I31").Value = 50002960
$ws.Range("K31").Value = 50002960
$ws.Range("M31").Value = -50002665

$ws.Range("H34").Value = 27030068
$ws.Range("I34").Value = 50002960
$ws.Range("K34").Value = 50002960
$ws.Range("M34").Value = -50002758

$ws.Range("H102").Value = 90332.664
$ws.Range("I102").Value = 42000
$ws.Range("J102").Value = 114499
$ws.Range("K102").Value = 42000
$ws.Range("L102").Value = 114499
$ws.Range("M102").Value = -39566
$ws.Range("N102").Value = -119367

$ws.Range("H122").Value = 5327.0713
$ws.Range("I122").Value = 5461.9
$ws.Range("J122").Value = 4990
$ws.Range("K122").Value = 16385.7
$ws.Range("L122").Value = 14970
$ws.Range("M122").Value = -13935.7
$ws.Range("N122").Value = -19870

$ws.Range("H132").Value = 2948.6924
$ws.Range("I132").Value = 2791.85
$ws.Range("J132").Value = 3471.5
$ws.Range("K132").Value = 8375.549999999999
$ws.Range("L132").Value = 10414.5
$ws.Range("M132").Value = -5845.549999999999
$ws.Range("N132").Value = -15474.5

$ws.Range("H134").Value = 2466.48
$ws.Range("I134").Value = 2402.5833
$ws.Range("K134").Value = 7207.749899999999
$ws.Range("M134").Value = -4672.749899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 823.1667
$ws.Range("I86").Value = 86.666664
$ws.Range("K86").Value = 259.999992
$ws.Range("M86").Value = 926.000008

$ws.Range("H89").Value = 823.1667
$ws.Range("I89").Value = 86.666664
$ws.Range("K89").Value = 779.9999759999999
$ws.Range("M89").Value = 5148.000024

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1938.6666
$ws.Range("J80").Value = 2573.8
$ws.Range("L80").Value = 2573.8
$ws.Range("N80").Value = -4569.8

$ws.Range("H83").Value = 1938.6666
$ws.Range("J83").Value = 2573.8
$ws.Range("L83").Value = 12869
$ws.Range("N83").Value = -22853

$ws.Range("H102").Value = 4002.3635
$ws.Range("I102").Value = 4150.5
$ws.Range("J102").Value = 3607.3333
$ws.Range("K102").Value = 4150.5
$ws.Range("L102").Value = 3607.3333
$ws.Range("M102").Value = -2528.5
$ws.Range("N102").Value = -6851.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 130899
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 130899
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 130899
$ws.Range("N6").Value = -131123
$ws.Range("M6").ClearContents()

$ws.Range("H46").Value = 3449.5
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H53").Value = 81001
$ws.Range("I53").Value = 81001
$ws.Range("K53").Value = 81001
$ws.Range("M53").Value = -80483

$ws.Range("H55").Value = 861.3200000000001
$ws.Range("I55").Value = 449.1
$ws.Range("J55").Value = 1136.1333
$ws.Range("K55").Value = 449.1
$ws.Range("L55").Value = 1136.1333
$ws.Range("M55").Value = -276.1
$ws.Range("N55").Value = -1482.1333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 11576.5
$ws.Range("I45").Value = 30484.5
$ws.Range("J45").Value = 5273.8335
$ws.Range("K45").Value = 30484.5
$ws.Range("L45").Value = 5273.8335
$ws.Range("M45").Value = -29993.5
$ws.Range("N45").Value = -6255.8335

$ws.Range("H55").Value = 29273
$ws.Range("I55").Value = 29546.5
$ws.Range("J55").Value = 28999.5
$ws.Range("K55").Value = 29546.5
$ws.Range("L55").Value = 28999.5
$ws.Range("M55").Value = -29269.5
$ws.Range("N55").Value = -29553.5

$ws.Range("H81").Value = 1606.3334
$ws.Range("I81").Value = 1606.3334
$ws.Range("K81").Value = 3212.6668
$ws.Range("M81").Value = -2151.6668

$ws.Range("H84").Value = 1606.3334
$ws.Range("I84").Value = 1606.3334
$ws.Range("K84").Value = 16063.334
$ws.Range("M84").Value = -10759.334

$ws.Range("H136").Value = 346025.84
$ws.Range("I136").Value = 13164.462
$ws.Range("J136").Value = 2509624.8
$ws.Range("K136").Value = 39493.386
$ws.Range("L136").Value = 7528874.399999999
$ws.Range("M136").Value = -36943.386
$ws.Range("N136").Value = -7533974.399999999
